$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.5694

$ws.Range("B4").Value = 5.953099999999997
$ws.Range("D4").Value = -8.035200000000001
$ws.Range("E4").Value = 13.3829

$ws.Range("D5").Value = -8.606500000000002

$ws.Range("B6").Value = 9.358499999999998

$ws.Range("B7").Value = 6.788499999999996

$ws.Range("D8").Value = -8.147400000000001

$ws.Range("E9").Value = 14.05670000000001

$ws.Range("E11").Value = 13.3551

$ws.Range("E14").Value = 13.03300000000001

$ws.Range("B16").Value = 8.274000000000008
$ws.Range("D16").Value = -7.929900000000004

$ws.Range("E18").Value = 12.5098

$ws.Range("B20").Value = 6.067699999999998

$ws.Range("D22").Value = -7.976200000000005

$ws.Range("E25").Value = 13.06409999999999
